$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.832.17"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.636.19"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'215.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "'0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.0642"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'19.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "1.638.69"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "1.861.76"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "0.0₃0769"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").Value = "'63.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "25.838.56"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'194.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'4.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'1.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").Value = "'139.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("E27").Value = "  -5.00%  "
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("D29").Value = "'15.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").Value = "'0.0495"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").Value = "'3.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "'3.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("D35").Value = "'2.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "'0.904"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "'0.554"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").Value = "1.115.96"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Value = "'5.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "'99.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "0.0₆0111"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").Value = "'2.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.06%  "
$ws.Range("E48").Value = "  -5.34%  "
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").Value = "'7.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("E51").Value = "  +0.46%  "
